$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Porcentaje de fallos total" (D) and "Porcentaje de fallo" (E) columns to 0,
# and refresh the recalculated "Rango" (F), "Maximo" (G) and "Minimo" (H) values
# for each movement row, as produced by re-running the quality report generation
# (showing the range when clicking the difficulty level).

$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 22.810807399706093
$ws.Range("G2").Value = 62.810807399706093
$ws.Range("H2").Value = 40

$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 28.19861889156364
$ws.Range("G3").Value = 38.19861889156364
$ws.Range("H3").Value = 10

$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 56.378670739293071
$ws.Range("G4").Value = 122.28134913072228
$ws.Range("H4").Value = 65.902678391429205

$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 43.355339105518681
$ws.Range("G5").Value = 57.986070536640781
$ws.Range("H5").Value = 14.630731431122097

$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 16.43556202943499
$ws.Range("G6").Value = 21.43556202943499
$ws.Range("H6").Value = 5

$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 6.2731813929354026
$ws.Range("G7").Value = 29.851647229533484
$ws.Range("H7").Value = 23.578465836598081
